$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E
$ws.Range("E4").Value = "note"

# E5 is a plain (non-shared) formula referencing D5
$ws.Range("E5").Formula = "=D5/2"

# E6:E25 share the same relative formula pattern
$ws.Range("E6:E25").Formula = "=D6/2"
